$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 276889.2
$ws.Range("I6").Value = 434643.56
$ws.Range("J6").Value = 819
$ws.Range("K6").Value = 1303930.68
$ws.Range("L6").Value = 2457
$ws.Range("M6").Value = -1303818.68
$ws.Range("N6").Value = -2681
$ws.Range("H9").Value = 276.66666
$ws.Range("I9").Value = 298.75
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 298.75
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -129.75
$ws.Range("N9").Value = -438
$ws.Range("H12").Value = 25338.1
$ws.Range("I12").Value = 250.8
$ws.Range("J12").Value = 100600
$ws.Range("K12").Value = 250.8
$ws.Range("L12").Value = 100600
$ws.Range("M12").Value = -80.80000000000001
$ws.Range("N12").Value = -100940
$ws.Range("H21").Value = 12215.625
$ws.Range("I21").Value = 12500
$ws.Range("J21").Value = 12175
$ws.Range("K21").Value = 12500
$ws.Range("L21").Value = 12175
$ws.Range("M21").Value = -12032
$ws.Range("N21").Value = -13111
$ws.Range("H23").Value = 12215.625
$ws.Range("I23").Value = 12500
$ws.Range("J23").Value = 12175
$ws.Range("K23").Value = 12500
$ws.Range("L23").Value = 12175
$ws.Range("M23").Value = -12266
$ws.Range("N23").Value = -12643
$ws.Range("H38").Value = 1112792.4
$ws.Range("I38").Value = 1536173.2
$ws.Range("J38").Value = 1417.5
$ws.Range("K38").Value = 4608519.6
$ws.Range("L38").Value = 4252.5
$ws.Range("M38").Value = -4608147.6
$ws.Range("N38").Value = -4996.5
$ws.Range("H39").Value = 882016.4399999999
$ws.Range("I39").Value = 1221140.5
$ws.Range("J39").Value = 294
$ws.Range("K39").Value = 3663421.5
$ws.Range("L39").Value = 882
$ws.Range("M39").Value = -3663125.5
$ws.Range("N39").Value = -1474
$ws.Range("H43").Value = 1523.75
$ws.Range("I43").Value = 1248.3334
$ws.Range("J43").Value = 1799.1666
$ws.Range("K43").Value = 1248.3334
$ws.Range("L43").Value = 1799.1666
$ws.Range("M43").Value = -1179.3334
$ws.Range("N43").Value = -1937.1666
$ws.Range("H58").Value = 1026952.44
$ws.Range("I58").Value = 3279024.5
$ws.Range("J58").Value = 3283.3635
$ws.Range("K58").Value = 9837073.5
$ws.Range("L58").Value = 9850.0905
$ws.Range("M58").Value = -9836923.5
$ws.Range("N58").Value = -10150.0905
$ws.Range("H98").Value = 1060.7142
$ws.Range("I98").Value = 1086.6666
$ws.Range("J98").Value = 905
$ws.Range("K98").Value = 1086.6666
$ws.Range("L98").Value = 905
$ws.Range("M98").Value = 411.3334
$ws.Range("N98").Value = -3901
$ws.Range("H122").Value = 1060.7142
$ws.Range("I122").Value = 1086.6666
$ws.Range("J122").Value = 905
$ws.Range("K122").Value = 3259.9998
$ws.Range("L122").Value = 2715
$ws.Range("M122").Value = -809.9998000000001
$ws.Range("N122").Value = -7615
$ws.Range("H132").Value = 10879223
$ws.Range("I132").Value = 12510336
$ws.Range("J132").Value = 5137
$ws.Range("K132").Value = 37531008
$ws.Range("L132").Value = 15411
$ws.Range("M132").Value = -37528478
$ws.Range("N132").Value = -20471
$ws.Range("H137").Value = 1625.2162
$ws.Range("I137").Value = 1302.5807
$ws.Range("J137").Value = 3292.1667
$ws.Range("K137").Value = 3907.7421
$ws.Range("L137").Value = 9876.500100000001
$ws.Range("M137").Value = -1357.7421
$ws.Range("N137").Value = -14976.5001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value = 18611
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 18611
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 18611
$ws.Range("N98").Value = -24601
$ws.Range("H122").Value = 2001.421
$ws.Range("I122").Value = 1988.6428
$ws.Range("J122").Value = 2037.2
$ws.Range("K122").Value = 5965.928400000001
$ws.Range("L122").Value = 6111.6
$ws.Range("M122").Value = -3515.928400000001
$ws.Range("N122").Value = -11011.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 18224.2
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 18224.2
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 18224.2
$ws.Range("N96").Value = -23716.2
$ws.Range("H132").Value = 2807.3147
$ws.Range("I132").Value = 2409.1462
$ws.Range("J132").Value = 4063.077
$ws.Range("K132").Value = 7227.4386
$ws.Range("L132").Value = 12189.231
$ws.Range("M132").Value = -4697.4386
$ws.Range("N132").Value = -17249.231

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1550.5
$ws.Range("I17").Value = 99
$ws.Range("J17").Value = 3002
$ws.Range("K17").Value = 297
$ws.Range("L17").Value = 9006
$ws.Range("M17").Value = -128
$ws.Range("N17").Value = -9344
$ws.Range("H34").Value = 477.3
$ws.Range("I34").Value = 114
$ws.Range("J34").Value = 840.6
$ws.Range("K34").Value = 342
$ws.Range("L34").Value = 2521.8
$ws.Range("M34").Value = -258
$ws.Range("N34").Value = -2689.8
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 8245.933999999999
$ws.Range("I55").Value = 20378
$ws.Range("J55").Value = 5819.52
$ws.Range("K55").Value = 61134
$ws.Range("L55").Value = 17458.56
$ws.Range("M55").Value = -60957
$ws.Range("N55").Value = -17812.56

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2782.1365
$ws.Range("I102").Value = 2794.8572
$ws.Range("J102").Value = 2759.875
$ws.Range("K102").Value = 2794.8572
$ws.Range("L102").Value = 2794.8572
$ws.Range("M102").Value = -1172.8572
$ws.Range("N102").Value = -6003.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 815.72
$ws.Range("I22").Value = 542.375
$ws.Range("J22").Value = 944.35297
$ws.Range("K22").Value = 542.375
$ws.Range("L22").Value = 944.35297
$ws.Range("M22").Value = -247.375
$ws.Range("N22").Value = -1534.35297
$ws.Range("H27").Value = 815.72
$ws.Range("I27").Value = 542.375
$ws.Range("J27").Value = 944.35297
$ws.Range("K27").Value = 542.375
$ws.Range("L27").Value = 944.35297
$ws.Range("M27").Value = -435.375
$ws.Range("N27").Value = -1158.35297
$ws.Range("H40").Value = 102089.8
$ws.Range("I40").Value = 501000
$ws.Range("J40").Value = 2362.25
$ws.Range("K40").Value = 501000
$ws.Range("L40").Value = 2362.25
$ws.Range("M40").Value = -500864
$ws.Range("N40").Value = -2634.25
$ws.Range("H46").Value = 723743.7
$ws.Range("I46").Value = 337.5
$ws.Range("J46").Value = 1013106.2
$ws.Range("K46").Value = 337.5
$ws.Range("L46").Value = 1013106.2
$ws.Range("M46").Value = -149.5
$ws.Range("N46").Value = -1013482.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 44464.125
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 44464.125
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 44464.125
$ws.Range("N135").Value = -54604.125
